$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format only on the specific cells whose new values would
# otherwise be auto-converted to numbers by Excel, losing formatting
# such as trailing zeros (e.g. "1.000", "0.4490", "34.70").
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range('B12').Value = 'TRON'
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('B23').Value = 'BitDAO'
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('B26').Value = 'Monero'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('B31').Value = 'Stellar'
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('B34').Value = 'Hedera'
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('B38').Value = 'VeChain'
$ws.Range('B39').Value = 'MXToken'
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('B41').Value = 'Aave'
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('B44').Value = 'Quant'
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('B47').Value = 'Aptos'
$ws.Range('B48').Value = 'Algorand'
$ws.Range('B50').Value = 'Elrond'
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D2').Value = '30.471.64'
$ws.Range('D3').Value = '1.913.22'
$ws.Range('D5').Value = '239.16'
$ws.Range('D8').Value = '0.2841'
$ws.Range('D9').Value = '0.06695'
$ws.Range('D10').Value = '18.81'
$ws.Range('D11').Value = '101.73'
$ws.Range('D12').Value = '0.07694'
$ws.Range('D13').Value = '1.913.39'
$ws.Range('D14').Value = '5.204'
$ws.Range('D15').Value = '0.6713'
$ws.Range('D16').Value = '265.74'
$ws.Range('D17').Value = '30.484.49'
$ws.Range('D19').Value = '0.000007465'
$ws.Range('D20').Value = '12.66'
$ws.Range('D21').Value = '5.394'
$ws.Range('D22').Value = '1.000'
$ws.Range('D23').Value = '0.4490'
$ws.Range('D24').Value = '6.296'
$ws.Range('D25').Value = '9.354'
$ws.Range('D26').Value = '167.19'
$ws.Range('D27').Value = '19.13'
$ws.Range('D28').Value = '2.060'
$ws.Range('D29').Value = '4.721'
$ws.Range('D30').Value = '1.384'
$ws.Range('D31').Value = '0.09982'
$ws.Range('D32').Value = '1.509'
$ws.Range('D33').Value = '4.248'
$ws.Range('D34').Value = '0.04713'
$ws.Range('D35').Value = '0.7263'
$ws.Range('D36').Value = '1.107'
$ws.Range('D37').Value = '2.713'
$ws.Range('D38').Value = '0.01915'
$ws.Range('D39').Value = '2.624'
$ws.Range('D40').Value = '6.233'
$ws.Range('D41').Value = '74.92'
$ws.Range('D42').Value = '1.970'
$ws.Range('D43').Value = '0.8623'
$ws.Range('D44').Value = '105.97'
$ws.Range('D45').Value = '0.4253'
$ws.Range('D46').Value = '1.000'
$ws.Range('D47').Value = '7.397'
$ws.Range('D48').Value = '0.1202'
$ws.Range('D49').Value = '917.42'
$ws.Range('D50').Value = '34.70'
$ws.Range('D51').Value = '8.758'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  -2.33%  '
$ws.Range('E8').Value = '  -3.85%  '
$ws.Range('E9').Value = '  -2.96%  '
$ws.Range('E10').Value = '  -3.16%  '
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('E13').Value = '  -1.46%  '
$ws.Range('E14').Value = '  -2.59%  '
$ws.Range('E15').Value = '  -3.95%  '
$ws.Range('E16').Value = '  -3.81%  '
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('E20').Value = '  -3.46%  '
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  -12.46%  '
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('E25').Value = '  -3.76%  '
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('E28').Value = '  -4.26%  '
$ws.Range('E29').Value = '  +3.66%  '
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('E31').Value = '  -4.62%  '
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('E33').Value = '  -2.69%  '
$ws.Range('E34').Value = '  -2.77%  '
$ws.Range('E35').Value = '  -3.22%  '
$ws.Range('E36').Value = '  -4.34%  '
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('E38').Value = '  -3.76%  '
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('E40').Value = '  -3.64%  '
$ws.Range('E41').Value = '  -4.67%  '
$ws.Range('E42').Value = '  -5.95%  '
$ws.Range('E43').Value = '  -4.86%  '
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('E45').Value = '  -3.21%  '
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('E47').Value = '  -4.53%  '
$ws.Range('E48').Value = '  -3.27%  '
$ws.Range('E49').Value = '  -7.03%  '
$ws.Range('E50').Value = '  -3.51%  '
$ws.Range('E51').Value = '  -5.20%  '
